$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 561, shifting existing rows 561:669 down to 562:670
$ws.Rows.Item(561).Insert()

# Populate the newly inserted row 561 with values (mirrors the row below it,
# which retains most of the static values, updated with new measurement data)
$ws.Cells.Item(561, 1).Value = 5
$ws.Cells.Item(561, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(561, 3).Value = "Maule"
$ws.Cells.Item(561, 4).Value = 45209
$ws.Cells.Item(561, 4).NumberFormat = $ws.Cells.Item(562, 4).NumberFormat
$ws.Cells.Item(561, 5).Value = 7
$ws.Cells.Item(561, 6).Value = 100114014
$ws.Cells.Item(561, 7).Value = "Betarraga"
$ws.Cells.Item(561, 8).Value = "Sin especificar"
$ws.Cells.Item(561, 9).Value = "Primera"
$ws.Cells.Item(561, 10).Value = 5000
$ws.Cells.Item(561, 11).Value = 500
$ws.Cells.Item(561, 12).Value = 500
$ws.Cells.Item(561, 13).Value = 500
$ws.Cells.Item(561, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(561, 15).Value = "Región del Maule"
$ws.Cells.Item(561, 16).Value = 100
$ws.Cells.Item(561, 17).Value = 5
$ws.Cells.Item(561, 18).Value = "Hortaliza"
